# C5-PowerPoint.pptx — Fri, May 22, 2020  5:05:47 PM
#
# 1) Turn on "Embed fonts in the file" for the deck and make sure the
#    Limelight typeface (used on slide 8) is embedded, mirroring
#    PowerPoint's File > Options > Save > "Embed fonts in the file"
#    workflow (p:presentation/@embedTrueTypeFonts + p:embeddedFontLst).
# 2) Re-apply the table style on the slide 6 table (new theme table
#    style GUID).

$p = $ppt.ActivePresentation

# --- 1. Embed TrueType fonts (Limelight) -------------------------------
$p.EmbedTrueTypeFonts = $true

$fonts = $p.Fonts
$limelight = $null
for ($i = 1; $i -le $fonts.Count; $i++) {
    $candidate = $fonts.Item($i)
    if ($candidate.NameOther -eq "Limelight") {
        $limelight = $candidate
    }
}
if (-not $limelight) {
    $limelight = $fonts.Add("Limelight")
}
$limelight.Embedded = $true

# --- 2. Swap the table style on slide 6's table -------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shape = $tableSlide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{A351C970-23AB-4AE1-919D-B6B908D818C6}")
    }
}
